$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 0.63699999032322
$ws.Range("F2").Value = 0.1999999066667279
$ws.Range("G2").Value = 0.05999988438146938

# Row 3
$ws.Range("E3").Value = 0.2580000136364074
$ws.Range("F3").Value = 0.6000001199999307
$ws.Range("G3").Value = 0.2309999275684518

# Row 4
$ws.Range("E4").Value = 0.1049999960403727
$ws.Range("F4").Value = 0.1999999733333413
$ws.Range("G4").Value = 0.7090001880500788

# Row 5
$ws.Range("C5").Value = 0.07200008696040985
$ws.Range("D5").Value = 0.05800016977734325

# Row 6
$ws.Range("C6").Value = 0.6490000112262161
$ws.Range("D6").Value = 0.2070004449070189

# Row 7
$ws.Range("C7").Value = 0.2789999018133741
$ws.Range("D7").Value = 0.7349993853156379
